$d = $word.ActiveDocument

# The LINCS-specific Heading5 "LINCS" paragraphs (and the FirstParagraph
# content paragraph that immediately follows each of them) are being
# removed from the AT: Awareness and Training section, per the commit
# "recreated docs w/o LINCS". Everything else (AWS / CivicActions
# sub-sections, the AT-1..AT-4 headings, and the body text blocks) stays
# exactly where it is; Word will renumber bookmark ids automatically.

$lincsBookmarks = @("lincs", "lincs-1", "lincs-2", "lincs-3", "lincs-4")

foreach ($bmName in $lincsBookmarks) {
    if (-not $d.Bookmarks.Exists($bmName)) {
        continue
    }

    $bm = $d.Bookmarks.Item($bmName)
    $bmRange = $bm.Range
    $bmStart = $bmRange.Start

    # Locate the paragraph that actually contains the bookmark (the
    # "LINCS" Heading5 paragraph). We can't trust Range.Paragraphs here
    # because the bookmark is a zero/short-length range that sits right
    # on a paragraph boundary, so walk $d.Paragraphs explicitly and use a
    # half-open interval test.
    $targetPara = $null
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($bmStart -ge $p.Range.Start -and $bmStart -lt $p.Range.End) {
            $targetPara = $p
            break
        }
    }

    if ($null -eq $targetPara) {
        continue
    }

    # The paragraph immediately after the "LINCS" heading is the
    # corresponding FirstParagraph body text that belongs to it; both get
    # removed together.
    $nextPara = $targetPara.Next()

    $deleteStart = $targetPara.Range.Start
    $deleteEnd = if ($null -ne $nextPara) { $nextPara.Range.End } else { $targetPara.Range.End }

    $deleteRange = $d.Range($deleteStart, $deleteEnd)
    $deleteRange.Delete()
}
